$d = $word.ActiveDocument

# 1. "the Tornado is carrying" -> "more modern aircraft such as the Tornado are carrying"
$d.Content.Find.Execute(
    "the Tornado is carrying a heavy load in the Luftwaffe.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "more modern aircraft such as the Tornado are carrying a heavy load in the Luftwaffe.",
    2) | Out-Null

# 2. "which will be included here." -> "which are included here."
$d.Content.Find.Execute(
    "acquired 112 IDS variants which will be included here.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "acquired 112 IDS variants which are included here.",
    2) | Out-Null

# 3. " it was perfectly suited the intra-theater" -> " the Transall was perfectly suited the intra-theater"
$d.Content.Find.Execute(
    " it was perfectly suited the intra-theater",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    " the Transall was perfectly suited the intra-theater",
    2) | Out-Null
